# Update the "取得日時" (retrieved datetime) column (A) for rows 2-9 on the
# active sheet ("ランサーズ") to reflect the latest scrape timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-09-16 06:34:48"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
